# Update crypto price/volume data per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.163.67'
$ws.Cells.Item(2, 5).Value = '  -1.96%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.896.84'
$ws.Cells.Item(3, 5).Value = '  -2.66%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9980'
$ws.Cells.Item(4, 5).Value = '  -0.25%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '331.39'
$ws.Cells.Item(5, 5).Value = '  -3.12%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '0.9963'
$ws.Cells.Item(6, 5).Value = '  -0.37%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.4602'
$ws.Cells.Item(7, 5).Value = '  -3.70%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.4117'
$ws.Cells.Item(8, 5).Value = '  -0.92%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '47.82'
$ws.Cells.Item(9, 5).Value = '  -1.03%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.07984'
$ws.Cells.Item(10, 5).Value = '  -3.48%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '1.002'
$ws.Cells.Item(11, 5).Value = '  -3.98%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '22.09'
$ws.Cells.Item(12, 5).Value = '  -2.71%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.898.26'
$ws.Cells.Item(13, 5).Value = '  -2.20%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '5.943'
$ws.Cells.Item(14, 5).Value = '  -4.10%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '7.099'
$ws.Cells.Item(15, 5).Value = '  -4.40%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'Litecoin'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(16, 4).Value = '89.27'
$ws.Cells.Item(16, 5).Value = '  -3.43%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'BinanceUSD'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(17, 4).Value = '0.9997'
$ws.Cells.Item(17, 5).Value = '  -0.16%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value = '0.00001028'
$ws.Cells.Item(18, 5).Value = '  -3.46%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'TRON'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(19, 4).Value = '0.06565'
$ws.Cells.Item(19, 5).Value = '  -1.57%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '17.53'
$ws.Cells.Item(20, 5).Value = '  -2.97%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.0000'
$ws.Cells.Item(21, 5).Value = '  +0.01%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '29.049.90'
$ws.Cells.Item(22, 5).Value = '  -2.23%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '5.447'
$ws.Cells.Item(23, 5).Value = '  -2.95%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '11.46'
$ws.Cells.Item(24, 5).Value = '  +1.75%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '2.202'
$ws.Cells.Item(25, 5).Value = '  -3.50%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '2.133.68'
$ws.Cells.Item(26, 5).Value = '  -1.86%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '156.41'
$ws.Cells.Item(27, 5).Value = '  -2.68%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '19.66'
$ws.Cells.Item(28, 5).Value = '  -2.78%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '2.117'
$ws.Cells.Item(29, 5).Value = '  -3.74%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '5.613'
$ws.Cells.Item(30, 5).Value = '  -1.11%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '117.38'
$ws.Cells.Item(31, 5).Value = '  -4.17%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '1.041'
$ws.Cells.Item(32, 5).Value = '  +1.14%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '0.09359'
$ws.Cells.Item(33, 5).Value = '  -2.90%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '1.419'
$ws.Cells.Item(34, 5).Value = '  -4.18%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.510'
$ws.Cells.Item(35, 5).Value = '  -4.56%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '5.343'
$ws.Cells.Item(36, 5).Value = '  -2.69%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '0.06062'
$ws.Cells.Item(37, 5).Value = '  -4.05%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '0.02235'
$ws.Cells.Item(38, 5).Value = '  -4.31%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '8.362'
$ws.Cells.Item(39, 5).Value = '  -2.96%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '1.174'
$ws.Cells.Item(40, 5).Value = '  -2.04%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'TheSandbox'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(41, 4).Value = '0.5799'
$ws.Cells.Item(41, 5).Value = '  -5.26%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Frax'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(42, 4).Value = '0.9991'
$ws.Cells.Item(42, 5).Value = '  -0.08%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '0.1826'
$ws.Cells.Item(43, 5).Value = '  -4.07%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '10.13'
$ws.Cells.Item(44, 5).Value = '  -5.63%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '1.263'
$ws.Cells.Item(45, 5).Value = '  -0.69%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '0.07552'
$ws.Cells.Item(46, 5).Value = '  +1.91%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '2.291'
$ws.Cells.Item(47, 5).Value = '  -5.16%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).Value = '12.08'
$ws.Cells.Item(48, 5).Value = '  -4.32%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Decentraland'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(49, 4).Value = '0.5474'
$ws.Cells.Item(49, 5).Value = '  -4.27%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '1.911'
$ws.Cells.Item(50, 5).Value = '  -4.68%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '111.94'
$ws.Cells.Item(51, 5).Value = '  -1.88%  '
